$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Code Review 4 (column E) scores for week 3 contribution update
$ws.Range("E2").Value = 25
$ws.Range("E3").Value = 25
$ws.Range("E4").Value = 25
$ws.Range("E5").Value = 25

# Recalculate formulas (E7 SUM should now total 100)
$excel.Calculate()

# Update the active selection to E6 as in the final saved state
$ws.Range("E6").Select()

$wb.Save()
